# Updates the cryptos list prices (column D) and 1h volume-change percentages
# (column E) to the latest scraped values, matching the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is plain numeric text (e.g. "6.38") need to be
# force-typed as text first, otherwise Excel auto-converts them to a Number and we
# lose the fixed-decimal formatting (e.g. trailing zero in "592.87"/"0.990").
$textCells = @("D5","D6","D11","D12","D20","D21","D22","D25","D26","D29","D38","D40","D41","D42","D44","D46","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.321.20"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "3.765.38"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "592.87"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "165.80"
$ws.Range("E6").Value = "  -3.49%  "
$ws.Range("D7").Value = "3.764.04"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("D11").Value = "6.38"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("E13").Value = "  -4.20%  "
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").Value = "4.397.09"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "3.761.16"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "67.294.64"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("D20").Value = "6.95"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("D21").Value = "10.23"
$ws.Range("E21").Value = "  -8.27%  "
$ws.Range("D22").Value = "456.40"
$ws.Range("E22").Value = "  -3.68%  "
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").Value = "83.06"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").Value = "11.84"
$ws.Range("E26").Value = "  -3.56%  "
$ws.Range("E27").Value = "  -6.50%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "9.97"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("E31").Value = "  -2.50%  "
$ws.Range("E32").Value = "  -4.44%  "
$ws.Range("E33").Value = "  -3.72%  "
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "3.717.29"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  -7.94%  "
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("D40").Value = "0.990"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").Value = "5.72"
$ws.Range("E41").Value = "  -3.13%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D44").Value = "43.49"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").Value = "46.82"
$ws.Range("E47").Value = "  -3.77%  "
$ws.Range("D48").Value = "147.45"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("E49").Value = "  -8.55%  "
$ws.Range("D50").Value = "388.67"
$ws.Range("E50").Value = "  -3.72%  "
$ws.Range("D51").Value = "2.742.54"
$ws.Range("E51").Value = "  +1.51%  "

# Restore the default (unstyled) cell style now that the text is safely stored,
# so the cells keep matching their original (style-less) formatting.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
